$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update RF (column I) values for rows 26 through 41 from the old
# computed value (14.04425531914894) to the new 2025 value (15.008)
for ($r = 26; $r -le 41; $r++) {
    $ws.Cells.Item($r, 9).Value = 15.008
}
